# Scheduled market-data refresh: updates computed price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) on the per-job Leve
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect newly
# scraped market-board averages. Columns H-N hold plain numeric values
# (no formulas), so each touched cell is just reassigned to its new
# snapshot value. A couple of rows gain/lose an M (LeveProfitNQ) cell
# entirely because that computed value only exists when its inputs are
# non-zero.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 32 (ALC)
$ws_ALC.Cells.Item(32, 8).Value2 = 3857.818  # H32: 3619.4583 -> 3857.818
$ws_ALC.Cells.Item(32, 9).Value2 = 3104.5  # I32: 2803.5715 -> 3104.5
$ws_ALC.Cells.Item(32, 10).Value2 = 4140.3125  # J32: 3955.4119 -> 4140.3125
$ws_ALC.Cells.Item(32, 11).Value2 = 3104.5  # K32: 2803.5715 -> 3104.5
$ws_ALC.Cells.Item(32, 12).Value2 = 4140.3125  # L32: 3955.4119 -> 4140.3125
$ws_ALC.Cells.Item(32, 13).Value2 = -2778.5  # M32: -2477.5715 -> -2778.5
$ws_ALC.Cells.Item(32, 14).Value2 = -4792.3125  # N32: -4607.4119 -> -4792.3125

# Row 82 (ALC)
$ws_ALC.Cells.Item(82, 8).Value2 = 1517.4  # H82: 1658.8889 -> 1517.4
$ws_ALC.Cells.Item(82, 9).Value2 = 1517.4  # I82: 1658.8889 -> 1517.4
$ws_ALC.Cells.Item(82, 11).Value2 = 4552.200000000001  # K82: 4976.6667 -> 4552.200000000001
$ws_ALC.Cells.Item(82, 13).Value2 = -4146.200000000001  # M82: -4570.6667 -> -4146.200000000001

# Row 85 (ALC)
$ws_ALC.Cells.Item(85, 8).Value2 = 1517.4  # H85: 1658.8889 -> 1517.4
$ws_ALC.Cells.Item(85, 9).Value2 = 1517.4  # I85: 1658.8889 -> 1517.4
$ws_ALC.Cells.Item(85, 11).Value2 = 4552.200000000001  # K85: 4976.6667 -> 4552.200000000001
$ws_ALC.Cells.Item(85, 13).Value2 = -3148.200000000001  # M85: -3572.6667 -> -3148.200000000001

# Row 112 (ALC)
$ws_ALC.Cells.Item(112, 8).Value2 = 1483.9736  # H112: 1423.6666 -> 1483.9736
$ws_ALC.Cells.Item(112, 9).Value2 = 500  # I112: 0 -> 500
$ws_ALC.Cells.Item(112, 10).Value2 = 1510.5676  # J112: 1423.6666 -> 1510.5676
$ws_ALC.Cells.Item(112, 11).Value2 = 1500  # K112: 0 -> 1500
$ws_ALC.Cells.Item(112, 12).Value2 = 4531.7028  # L112: 4270.9998 -> 4531.7028
$ws_ALC.Cells.Item(112, 13).Value2 = -392  # M112: None -> -392
$ws_ALC.Cells.Item(112, 14).Value2 = -6747.7028  # N112: -6486.9998 -> -6747.7028

# Row 113 (ALC)
$ws_ALC.Cells.Item(113, 8).Value2 = 5659.8784  # H113: 5593.72 -> 5659.8784
$ws_ALC.Cells.Item(113, 9).Value2 = 6148.933  # I113: 6481 -> 6148.933
$ws_ALC.Cells.Item(113, 10).Value2 = 5535.5425  # J113: 5390.082 -> 5535.5425
$ws_ALC.Cells.Item(113, 11).Value2 = 6148.933  # K113: 6481 -> 6148.933
$ws_ALC.Cells.Item(113, 12).Value2 = 5535.5425  # L113: 5390.082 -> 5535.5425
$ws_ALC.Cells.Item(113, 13).Value2 = -2894.933  # M113: -3227 -> -2894.933
$ws_ALC.Cells.Item(113, 14).Value2 = -12043.5425  # N113: -11898.082 -> -12043.5425

# Row 116 (ALC)
$ws_ALC.Cells.Item(116, 8).Value2 = 3647  # H116: 3461.8386 -> 3647
$ws_ALC.Cells.Item(116, 9).Value2 = 3232.8462  # I116: 2951.0625 -> 3232.8462
$ws_ALC.Cells.Item(116, 10).Value2 = 4005.9333  # J116: 4006.6667 -> 4005.9333
$ws_ALC.Cells.Item(116, 11).Value2 = 3232.8462  # K116: 2951.0625 -> 3232.8462
$ws_ALC.Cells.Item(116, 12).Value2 = 4005.9333  # L116: 4006.6667 -> 4005.9333
$ws_ALC.Cells.Item(116, 13).Value2 = 209.1538  # M116: 490.9375 -> 209.1538
$ws_ALC.Cells.Item(116, 14).Value2 = -10889.9333  # N116: -10890.6667 -> -10889.9333

# Row 138 (ALC)
$ws_ALC.Cells.Item(138, 8).Value2 = 1769.94  # H138: 1705.66 -> 1769.94
$ws_ALC.Cells.Item(138, 9).Value2 = 907.1739  # I138: 877.5417 -> 907.1739
$ws_ALC.Cells.Item(138, 10).Value2 = 2504.889  # J138: 2470.077 -> 2504.889
$ws_ALC.Cells.Item(138, 11).Value2 = 2721.5217  # K138: 2632.6251 -> 2721.5217
$ws_ALC.Cells.Item(138, 12).Value2 = 7514.667  # L138: 7410.231000000001 -> 7514.667
$ws_ALC.Cells.Item(138, 13).Value2 = 2418.4783  # M138: 2507.3749 -> 2418.4783
$ws_ALC.Cells.Item(138, 14).Value2 = -17794.667  # N138: -17690.231 -> -17794.667

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws_ARM.Cells.Item(2, 8).Value2 = 5317.4287  # H2: 5660.3076 -> 5317.4287
$ws_ARM.Cells.Item(2, 9).Value2 = 1553.619  # I2: 1626.6316 -> 1553.619
$ws_ARM.Cells.Item(2, 11).Value2 = 1553.619  # K2: 1626.6316 -> 1553.619
$ws_ARM.Cells.Item(2, 13).Value2 = -1440.619  # M2: -1513.6316 -> -1440.619

# Row 6 (ARM)
$ws_ARM.Cells.Item(6, 8).Value2 = 0  # H6: 2995 -> 0
$ws_ARM.Cells.Item(6, 9).Value2 = 0  # I6: 2995 -> 0
$ws_ARM.Cells.Item(6, 11).Value2 = 0  # K6: 2995 -> 0
$ws_ARM.Cells.Item(6, 13).ClearContents()  # M6 removed (was -2822)

# Row 45 (ARM)
$ws_ARM.Cells.Item(45, 8).Value2 = 1118.3158  # H45: 1127.3334 -> 1118.3158
$ws_ARM.Cells.Item(45, 9).Value2 = 1132.5385  # I45: 1147.25 -> 1132.5385
$ws_ARM.Cells.Item(45, 11).Value2 = 1132.5385  # K45: 1147.25 -> 1132.5385
$ws_ARM.Cells.Item(45, 13).Value2 = -755.5385000000001  # M45: -770.25 -> -755.5385000000001

# Row 116 (ARM)
$ws_ARM.Cells.Item(116, 8).Value2 = 5317.4287  # H116: 5660.3076 -> 5317.4287
$ws_ARM.Cells.Item(116, 9).Value2 = 1553.619  # I116: 1626.6316 -> 1553.619
$ws_ARM.Cells.Item(116, 11).Value2 = 1553.619  # K116: 1626.6316 -> 1553.619
$ws_ARM.Cells.Item(116, 13).Value2 = 740.3810000000001  # M116: 667.3684000000001 -> 740.3810000000001

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws_BSM.Cells.Item(3, 8).Value2 = 5317.4287  # H3: 5660.3076 -> 5317.4287
$ws_BSM.Cells.Item(3, 9).Value2 = 1553.619  # I3: 1626.6316 -> 1553.619
$ws_BSM.Cells.Item(3, 11).Value2 = 1553.619  # K3: 1626.6316 -> 1553.619
$ws_BSM.Cells.Item(3, 13).Value2 = -1439.619  # M3: -1512.6316 -> -1439.619

# Row 86 (BSM)
$ws_BSM.Cells.Item(86, 8).Value2 = 2809.25  # H86: 2312.5938 -> 2809.25
$ws_BSM.Cells.Item(86, 9).Value2 = 2257  # I86: 1756.25 -> 2257
$ws_BSM.Cells.Item(86, 10).Value2 = 3461.9092  # J86: 3239.8333 -> 3461.9092
$ws_BSM.Cells.Item(86, 11).Value2 = 2257  # K86: 1756.25 -> 2257
$ws_BSM.Cells.Item(86, 12).Value2 = 3461.9092  # L86: 3239.8333 -> 3461.9092
$ws_BSM.Cells.Item(86, 13).Value2 = -1134  # M86: -633.25 -> -1134
$ws_BSM.Cells.Item(86, 14).Value2 = -5707.9092  # N86: -5485.8333 -> -5707.9092

# Row 89 (BSM)
$ws_BSM.Cells.Item(89, 8).Value2 = 2809.25  # H89: 2312.5938 -> 2809.25
$ws_BSM.Cells.Item(89, 9).Value2 = 2257  # I89: 1756.25 -> 2257
$ws_BSM.Cells.Item(89, 10).Value2 = 3461.9092  # J89: 3239.8333 -> 3461.9092
$ws_BSM.Cells.Item(89, 11).Value2 = 11285  # K89: 8781.25 -> 11285
$ws_BSM.Cells.Item(89, 12).Value2 = 17309.546  # L89: 16199.1665 -> 17309.546
$ws_BSM.Cells.Item(89, 13).Value2 = -5669  # M89: -3165.25 -> -5669
$ws_BSM.Cells.Item(89, 14).Value2 = -28541.546  # N89: -27431.1665 -> -28541.546

# Row 107 (BSM)
$ws_BSM.Cells.Item(107, 8).Value2 = 2391.875  # H107: 2444.842 -> 2391.875
$ws_BSM.Cells.Item(107, 9).Value2 = 1319.7693  # I107: 1556 -> 1319.7693
$ws_BSM.Cells.Item(107, 10).Value2 = 7037.6665  # J107: 10000 -> 7037.6665
$ws_BSM.Cells.Item(107, 11).Value2 = 1319.7693  # K107: 1556 -> 1319.7693
$ws_BSM.Cells.Item(107, 12).Value2 = 7037.6665  # L107: 10000 -> 7037.6665
$ws_BSM.Cells.Item(107, 13).Value2 = 600.2307000000001  # M107: 364 -> 600.2307000000001
$ws_BSM.Cells.Item(107, 14).Value2 = -10877.6665  # N107: -13840 -> -10877.6665

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws_CRP.Cells.Item(16, 8).Value2 = 3253.6428  # H16: 3177.1614 -> 3253.6428
$ws_CRP.Cells.Item(16, 9).Value2 = 3672.8333  # I16: 3538.4443 -> 3672.8333
$ws_CRP.Cells.Item(16, 11).Value2 = 3672.8333  # K16: 3538.4443 -> 3672.8333
$ws_CRP.Cells.Item(16, 13).Value2 = -3385.8333  # M16: -3251.4443 -> -3385.8333

# Row 31 (CRP)
$ws_CRP.Cells.Item(31, 8).Value2 = 1577.625  # H31: 1491.4117 -> 1577.625
$ws_CRP.Cells.Item(31, 9).Value2 = 1865.1  # I31: 1705.7273 -> 1865.1
$ws_CRP.Cells.Item(31, 11).Value2 = 1865.1  # K31: 1705.7273 -> 1865.1
$ws_CRP.Cells.Item(31, 13).Value2 = -1570.1  # M31: -1410.7273 -> -1570.1

# Row 34 (CRP)
$ws_CRP.Cells.Item(34, 8).Value2 = 1577.625  # H34: 1491.4117 -> 1577.625
$ws_CRP.Cells.Item(34, 9).Value2 = 1865.1  # I34: 1705.7273 -> 1865.1
$ws_CRP.Cells.Item(34, 11).Value2 = 1865.1  # K34: 1705.7273 -> 1865.1
$ws_CRP.Cells.Item(34, 13).Value2 = -1663.1  # M34: -1503.7273 -> -1663.1

# Row 113 (CRP)
$ws_CRP.Cells.Item(113, 8).Value2 = 3253.6428  # H113: 3177.1614 -> 3253.6428
$ws_CRP.Cells.Item(113, 9).Value2 = 3672.8333  # I113: 3538.4443 -> 3672.8333
$ws_CRP.Cells.Item(113, 11).Value2 = 3672.8333  # K113: 3538.4443 -> 3672.8333
$ws_CRP.Cells.Item(113, 13).Value2 = -1502.8333  # M113: -1368.4443 -> -1502.8333

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws_CUL.Cells.Item(5, 8).Value2 = 1103.2354  # H5: 896.56525 -> 1103.2354
$ws_CUL.Cells.Item(5, 9).Value2 = 987.3333  # I5: 816.3 -> 987.3333
$ws_CUL.Cells.Item(5, 10).Value2 = 1972.5  # J5: 1431.6666 -> 1972.5
$ws_CUL.Cells.Item(5, 11).Value2 = 2961.9999  # K5: 2448.9 -> 2961.9999
$ws_CUL.Cells.Item(5, 12).Value2 = 5917.5  # L5: 4294.9998 -> 5917.5
$ws_CUL.Cells.Item(5, 13).Value2 = -2849.9999  # M5: -2336.9 -> -2849.9999
$ws_CUL.Cells.Item(5, 14).Value2 = -6141.5  # N5: -4518.9998 -> -6141.5

# Row 38 (CUL)
$ws_CUL.Cells.Item(38, 8).Value2 = 149.5  # H38: 160.92308 -> 149.5
$ws_CUL.Cells.Item(38, 9).Value2 = 91.666664  # I38: 109.8 -> 91.666664
$ws_CUL.Cells.Item(38, 11).Value2 = 274.999992  # K38: 329.4 -> 274.999992
$ws_CUL.Cells.Item(38, 13).Value2 = 72.00000799999998  # M38: 17.60000000000002 -> 72.00000799999998

# Row 135 (CUL)
$ws_CUL.Cells.Item(135, 8).Value2 = 1103.2354  # H135: 896.56525 -> 1103.2354
$ws_CUL.Cells.Item(135, 9).Value2 = 987.3333  # I135: 816.3 -> 987.3333
$ws_CUL.Cells.Item(135, 10).Value2 = 1972.5  # J135: 1431.6666 -> 1972.5
$ws_CUL.Cells.Item(135, 11).Value2 = 8885.9997  # K135: 7346.7 -> 8885.9997
$ws_CUL.Cells.Item(135, 12).Value2 = 17752.5  # L135: 12884.9994 -> 17752.5
$ws_CUL.Cells.Item(135, 13).Value2 = -6350.9997  # M135: -4811.7 -> -6350.9997
$ws_CUL.Cells.Item(135, 14).Value2 = -22822.5  # N135: -17954.9994 -> -22822.5

# Row 139 (CUL)
$ws_CUL.Cells.Item(139, 8).Value2 = 4659.8184  # H139: 5762.1113 -> 4659.8184
$ws_CUL.Cells.Item(139, 10).Value2 = 4879.8  # J139: 8333.333000000001 -> 4879.8
$ws_CUL.Cells.Item(139, 12).Value2 = 14639.4  # L139: 24999.999 -> 14639.4
$ws_CUL.Cells.Item(139, 14).Value2 = -24919.4  # N139: -35279.999 -> -24919.4

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 3 (GSM)
$ws_GSM.Cells.Item(3, 8).Value2 = 3665.3333  # H3: 3998 -> 3665.3333
$ws_GSM.Cells.Item(3, 9).Value2 = 3665.3333  # I3: 3998 -> 3665.3333
$ws_GSM.Cells.Item(3, 11).Value2 = 3665.3333  # K3: 3998 -> 3665.3333
$ws_GSM.Cells.Item(3, 13).Value2 = -3549.3333  # M3: -3882 -> -3549.3333

# Row 107 (GSM)
$ws_GSM.Cells.Item(107, 8).Value2 = 871.4  # H107: 646.63635 -> 871.4
$ws_GSM.Cells.Item(107, 9).Value2 = 908.8  # I107: 890.6667 -> 908.8
$ws_GSM.Cells.Item(107, 10).Value2 = 834  # J107: 353.8 -> 834
$ws_GSM.Cells.Item(107, 11).Value2 = 908.8  # K107: 890.6667 -> 908.8
$ws_GSM.Cells.Item(107, 12).Value2 = 834  # L107: 353.8 -> 834
$ws_GSM.Cells.Item(107, 13).Value2 = 1011.2  # M107: 1029.3333 -> 1011.2
$ws_GSM.Cells.Item(107, 14).Value2 = -4674  # N107: -4193.8 -> -4674

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws_LTW.Cells.Item(46, 8).Value2 = 9260071  # H46: 8334161 -> 9260071
$ws_LTW.Cells.Item(46, 9).Value2 = 13889607  # I46: 11905516 -> 13889607
$ws_LTW.Cells.Item(46, 11).Value2 = 13889607  # K46: 11905516 -> 13889607
$ws_LTW.Cells.Item(46, 13).Value2 = -13889419  # M46: -11905328 -> -13889419

# Row 61 (LTW)
$ws_LTW.Cells.Item(61, 8).Value2 = 1851.5  # H61: 1773.909 -> 1851.5
$ws_LTW.Cells.Item(61, 9).Value2 = 1758.125  # I61: 1673.6666 -> 1758.125
$ws_LTW.Cells.Item(61, 11).Value2 = 1758.125  # K61: 1673.6666 -> 1758.125
$ws_LTW.Cells.Item(61, 13).Value2 = -1556.125  # M61: -1471.6666 -> -1556.125

# Row 113 (LTW)
$ws_LTW.Cells.Item(113, 8).Value2 = 1851.5  # H113: 1773.909 -> 1851.5
$ws_LTW.Cells.Item(113, 9).Value2 = 1758.125  # I113: 1673.6666 -> 1758.125
$ws_LTW.Cells.Item(113, 11).Value2 = 1758.125  # K113: 1673.6666 -> 1758.125
$ws_LTW.Cells.Item(113, 13).Value2 = 411.875  # M113: 496.3334 -> 411.875

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 10 (WVR)
$ws_WVR.Cells.Item(10, 8).Value2 = 41000  # H10: 0 -> 41000
$ws_WVR.Cells.Item(10, 9).Value2 = 41000  # I10: 0 -> 41000
$ws_WVR.Cells.Item(10, 11).Value2 = 41000  # K10: 0 -> 41000
$ws_WVR.Cells.Item(10, 13).Value2 = -40831  # M10: None -> -40831

# Row 41 (WVR)
$ws_WVR.Cells.Item(41, 8).Value2 = 21103.8  # H41: 20253 -> 21103.8
$ws_WVR.Cells.Item(41, 9).Value2 = 24989  # I41: 21632 -> 24989
$ws_WVR.Cells.Item(41, 10).Value2 = 20132.5  # J41: 19701.4 -> 20132.5
$ws_WVR.Cells.Item(41, 11).Value2 = 24989  # K41: 21632 -> 24989
$ws_WVR.Cells.Item(41, 12).Value2 = 20132.5  # L41: 19701.4 -> 20132.5
$ws_WVR.Cells.Item(41, 13).Value2 = -24599  # M41: -21242 -> -24599
$ws_WVR.Cells.Item(41, 14).Value2 = -20912.5  # N41: -20481.4 -> -20912.5

# Row 51 (WVR)
$ws_WVR.Cells.Item(51, 8).Value2 = 79999  # H51: 47499.5 -> 79999
$ws_WVR.Cells.Item(51, 10).Value2 = 79999  # J51: 47499.5 -> 79999
$ws_WVR.Cells.Item(51, 12).Value2 = 79999  # L51: 47499.5 -> 79999
$ws_WVR.Cells.Item(51, 14).Value2 = -81019  # N51: -48519.5 -> -81019

# Row 81 (WVR)
$ws_WVR.Cells.Item(81, 8).Value2 = 3803.9333  # H81: 4097.0713 -> 3803.9333
$ws_WVR.Cells.Item(81, 9).Value2 = 1704.091  # I81: 1904.5 -> 1704.091
$ws_WVR.Cells.Item(81, 11).Value2 = 3408.182  # K81: 3809 -> 3408.182
$ws_WVR.Cells.Item(81, 13).Value2 = -2347.182  # M81: -2748 -> -2347.182

# Row 84 (WVR)
$ws_WVR.Cells.Item(84, 8).Value2 = 3803.9333  # H84: 4097.0713 -> 3803.9333
$ws_WVR.Cells.Item(84, 9).Value2 = 1704.091  # I84: 1904.5 -> 1704.091
$ws_WVR.Cells.Item(84, 11).Value2 = 17040.91  # K84: 19045 -> 17040.91
$ws_WVR.Cells.Item(84, 13).Value2 = -11736.91  # M84: -13741 -> -11736.91

# Row 113 (WVR)
$ws_WVR.Cells.Item(113, 8).Value2 = 1499.5714  # H113: 1433.5946 -> 1499.5714
$ws_WVR.Cells.Item(113, 10).Value2 = 1984.0625  # J113: 1794.6111 -> 1984.0625
$ws_WVR.Cells.Item(113, 12).Value2 = 5952.1875  # L113: 5383.8333 -> 5952.1875
$ws_WVR.Cells.Item(113, 14).Value2 = -10292.1875  # N113: -9723.8333 -> -10292.1875
